$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A10 was stored as text ("71277628"); the redemption is now confirmed,
# so it gets recorded as a real number like the other phone-number cells.
$ws.Range("A10").Value = 71277628

# Append the new redemption row for phone 71277628 / 766 points.
# Force column A to text first so the phone number round-trips as a
# string (matching the existing phone-as-text convention), then restore
# the default cell style so no stray number formatting is left behind.
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "71277628"
$ws.Range("A11").Style = "Normal"

$ws.Range("B11").Value = 766
$ws.Range("C11").Value = "2025-08-18T16:54:29"
